$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "1" to "Sighnaghi"
$ws.Name = "Sighnaghi"

# Update table values: several numeric cells in the Urban/Rural rows become
# "confidential/unavailable" markers ("…" already exists as a shared string,
# "..." is a new one introduced by this edit).
$ws.Range("C6").Value = "…"
$ws.Range("F6").Value = "…"
$ws.Range("G6").Value = "…"
$ws.Range("L6").Value = "…"
$ws.Range("B6").Value = "..."
$ws.Range("N6").Value = "..."
$ws.Range("O6").Value = "..."

$ws.Range("G7").Value = "…"
$ws.Range("B7").Value = "..."
$ws.Range("M7").Value = "..."
$ws.Range("N7").Value = "..."
$ws.Range("O7").Value = "..."

# Remove the blank spacer row (old row 8) so the note row moves up from 9 to 8
$ws.Range("A8").EntireRow.Delete()
